$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.039369959750353
$ws.Cells.Item(2, 4).Value2 = 1.045329498178673
$ws.Cells.Item(2, 5).Value2 = 1.047346641722156
$ws.Cells.Item(2, 6).Value2 = 1.056857548769779
$ws.Cells.Item(2, 9).Value2 = 1.037451470814696
$ws.Cells.Item(2, 10).Value2 = 1.044462345111604
$ws.Cells.Item(2, 11).Value2 = 1.048098043204668
$ws.Cells.Item(2, 12).Value2 = 1.050109535099642
$ws.Cells.Item(2, 13).Value2 = 1.059594109347404
$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.040298039018768
$ws.Cells.Item(3, 4).Value2 = 1.046016249075006
$ws.Cells.Item(3, 5).Value2 = 1.048144585214739
$ws.Cells.Item(3, 6).Value2 = 1.05767439403293
$ws.Cells.Item(3, 9).Value2 = 1.037593281652592
$ws.Cells.Item(3, 10).Value2 = 1.04503570918404
$ws.Cells.Item(3, 11).Value2 = 1.048596537367496
$ws.Cells.Item(3, 12).Value2 = 1.05071934034515
$ws.Cells.Item(3, 13).Value2 = 1.060224668439941
$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.040899290802824
$ws.Cells.Item(4, 4).Value2 = 1.046461067549611
$ws.Cells.Item(4, 5).Value2 = 1.048661859337785
$ws.Cells.Item(4, 6).Value2 = 1.058203792220207
$ws.Cells.Item(4, 9).Value2 = 1.037683906595885
$ws.Cells.Item(4, 10).Value2 = 1.045406803599949
$ws.Cells.Item(4, 11).Value2 = 1.0489188619897
$ws.Cells.Item(4, 12).Value2 = 1.051114219638494
$ws.Cells.Item(4, 13).Value2 = 1.060632872886258
$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.041152228640804
$ws.Cells.Item(5, 4).Value2 = 1.046648174014058
$ws.Cells.Item(5, 5).Value2 = 1.048879547430179
$ws.Cells.Item(5, 6).Value2 = 1.058426551644131
$ws.Cells.Item(5, 9).Value2 = 1.037721732689797
$ws.Cells.Item(5, 10).Value2 = 1.045562831886992
$ws.Cells.Item(5, 11).Value2 = 1.049054309870928
$ws.Cells.Item(5, 12).Value2 = 1.051280296052599
$ws.Cells.Item(5, 13).Value2 = 1.060804526184541
$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.041194708039192
$ws.Cells.Item(6, 4).Value2 = 1.046679596107805
$ws.Cells.Item(6, 5).Value2 = 1.048916111425125
$ws.Cells.Item(6, 6).Value2 = 1.058463965632557
$ws.Cells.Item(6, 9).Value2 = 1.037728067851678
$ws.Cells.Item(6, 10).Value2 = 1.04558903087122
$ws.Cells.Item(6, 11).Value2 = 1.049077048764012
$ws.Cells.Item(6, 12).Value2 = 1.051308185029576
$ws.Cells.Item(6, 13).Value2 = 1.060833350074855
$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.040902669898073
$ws.Cells.Item(7, 4).Value2 = 1.046463567265768
$ws.Cells.Item(7, 5).Value2 = 1.048664767210638
$ws.Cells.Item(7, 6).Value2 = 1.058206767957247
$ws.Cells.Item(7, 9).Value2 = 1.037684413102054
$ws.Cells.Item(7, 10).Value2 = 1.045408888379757
$ws.Cells.Item(7, 11).Value2 = 1.048920672078011
$ws.Cells.Item(7, 12).Value2 = 1.051116438489615
$ws.Cells.Item(7, 13).Value2 = 1.060635166354401
$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.039683458076469
$ws.Cells.Item(8, 4).Value2 = 1.045561495761001
$ws.Cells.Item(8, 5).Value2 = 1.047616112461954
$ws.Cells.Item(8, 6).Value2 = 1.05713342947654
$ws.Cells.Item(8, 9).Value2 = 1.037499631316672
$ws.Cells.Item(8, 10).Value2 = 1.044656096904436
$ws.Cells.Item(8, 11).Value2 = 1.04826655942103
$ws.Cells.Item(8, 12).Value2 = 1.050315559842655
$ws.Cells.Item(8, 13).Value2 = 1.059807169407394
$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.037540644539925
$ws.Cells.Item(9, 4).Value2 = 1.043975411225252
$ws.Cells.Item(9, 5).Value2 = 1.045775609913113
$ws.Cells.Item(9, 6).Value2 = 1.055248618418613
$ws.Cells.Item(9, 9).Value2 = 1.037165348757322
$ws.Cells.Item(9, 10).Value2 = 1.043330322460074
$ws.Cells.Item(9, 11).Value2 = 1.047112189810762
$ws.Cells.Item(9, 12).Value2 = 1.048906624723262
$ws.Cells.Item(9, 13).Value2 = 1.058349657305174
$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.036115943787897
$ws.Cells.Item(10, 4).Value2 = 1.042920462350672
$ws.Cells.Item(10, 5).Value2 = 1.044553658494293
$ws.Cells.Item(10, 6).Value2 = 1.053996587244456
$ws.Cells.Item(10, 9).Value2 = 1.036936700029401
$ws.Cells.Item(10, 10).Value2 = 1.042447044078303
$ws.Cells.Item(10, 11).Value2 = 1.046341515791822
$ws.Cells.Item(10, 12).Value2 = 1.047968971671646
$ws.Cells.Item(10, 13).Value2 = 1.057379094814896
$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.035499960657339
$ws.Cells.Item(11, 4).Value2 = 1.042464257962953
$ws.Cells.Item(11, 5).Value2 = 1.044025759185667
$ws.Cells.Item(11, 6).Value2 = 1.053455535995426
$ws.Cells.Item(11, 9).Value2 = 1.036836326363778
$ws.Cells.Item(11, 10).Value2 = 1.042064725702048
$ws.Cells.Item(11, 11).Value2 = 1.046007561708655
$ws.Cells.Item(11, 12).Value2 = 1.047563362707433
$ws.Cells.Item(11, 13).Value2 = 1.056959112730127
$ws.Cells.Item(12, 2).Value2 = 1.02
$ws.Cells.Item(12, 3).Value2 = 1.035271296901633
$ws.Cells.Item(12, 4).Value2 = 1.042294894506254
$ws.Cells.Item(12, 5).Value2 = 1.04382985794328
$ws.Cells.Item(12, 6).Value2 = 1.053254730519477
$ws.Cells.Item(12, 9).Value2 = 1.036798838320858
$ws.Cells.Item(12, 10).Value2 = 1.041922738948623
$ws.Cells.Item(12, 11).Value2 = 1.045883480424825
$ws.Cells.Item(12, 12).Value2 = 1.047412763076285
$ws.Cells.Item(12, 13).Value2 = 1.056803155905701
$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.035320339705524
$ws.Cells.Item(13, 4).Value2 = 1.042331219389804
$ws.Cells.Item(13, 5).Value2 = 1.043871871067844
$ws.Cells.Item(13, 6).Value2 = 1.053297796477847
$ws.Cells.Item(13, 9).Value2 = 1.036806888890718
$ws.Cells.Item(13, 10).Value2 = 1.041953194514436
$ws.Cells.Item(13, 11).Value2 = 1.045910097884234
$ws.Cells.Item(13, 12).Value2 = 1.047445064399228
$ws.Cells.Item(13, 13).Value2 = 1.056836607194423
$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.035481056380249
$ws.Cells.Item(14, 4).Value2 = 1.042450256457171
$ws.Cells.Item(14, 5).Value2 = 1.044009562163045
$ws.Cells.Item(14, 6).Value2 = 1.053438933974729
$ws.Cells.Item(14, 9).Value2 = 1.036833231764603
$ws.Cells.Item(14, 10).Value2 = 1.042052988549994
$ws.Cells.Item(14, 11).Value2 = 1.045997305835417
$ws.Cells.Item(14, 12).Value2 = 1.047550912830537
$ws.Cells.Item(14, 13).Value2 = 1.05694622040014
$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.035580097793384
$ws.Cells.Item(15, 4).Value2 = 1.042523611254451
$ws.Cells.Item(15, 5).Value2 = 1.044094422633029
$ws.Cells.Item(15, 6).Value2 = 1.053525915369376
$ws.Cells.Item(15, 9).Value2 = 1.036849435359978
$ws.Cells.Item(15, 10).Value2 = 1.042114478075105
$ws.Cells.Item(15, 11).Value2 = 1.04605103281729
$ws.Cells.Item(15, 12).Value2 = 1.047616137745302
$ws.Cells.Item(15, 13).Value2 = 1.057013762485229
$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.036156844092171
$ws.Cells.Item(16, 4).Value2 = 1.042950751824387
$ws.Cells.Item(16, 5).Value2 = 1.044588719151138
$ws.Cells.Item(16, 6).Value2 = 1.054032518080191
$ws.Cells.Item(16, 9).Value2 = 1.036943332745412
$ws.Cells.Item(16, 10).Value2 = 1.04247242049175
$ws.Cells.Item(16, 11).Value2 = 1.046363674125168
$ws.Cells.Item(16, 12).Value2 = 1.04799589916674
$ws.Cells.Item(16, 13).Value2 = 1.057406973623033
$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.036518869510398
$ws.Cells.Item(17, 4).Value2 = 1.043218846440379
$ws.Cells.Item(17, 5).Value2 = 1.044899104457876
$ws.Cells.Item(17, 6).Value2 = 1.054350588745859
$ws.Cells.Item(17, 9).Value2 = 1.037001866478545
$ws.Cells.Item(17, 10).Value2 = 1.042696988535167
$ws.Cells.Item(17, 11).Value2 = 1.046559720480956
$ws.Cells.Item(17, 12).Value2 = 1.048234221785965
$ws.Cells.Item(17, 13).Value2 = 1.057653700017349
$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.036730121652672
$ws.Cells.Item(18, 4).Value2 = 1.043375278830633
$ws.Cells.Item(18, 5).Value2 = 1.045080263947129
$ws.Cells.Item(18, 6).Value2 = 1.054536218694915
$ws.Cells.Item(18, 9).Value2 = 1.037035876252537
$ws.Cells.Item(18, 10).Value2 = 1.04282798937551
$ws.Cells.Item(18, 11).Value2 = 1.046674047072593
$ws.Cells.Item(18, 12).Value2 = 1.048373269937081
$ws.Cells.Item(18, 13).Value2 = 1.057797638061267
$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.036802168202466
$ws.Cells.Item(19, 4).Value2 = 1.04342862790643
$ws.Cells.Item(19, 5).Value2 = 1.045142054442838
$ws.Cells.Item(19, 6).Value2 = 1.054599531416117
$ws.Cells.Item(19, 9).Value2 = 1.037047450304634
$ws.Cells.Item(19, 10).Value2 = 1.042872659630892
$ws.Cells.Item(19, 11).Value2 = 1.046713025388371
$ws.Cells.Item(19, 12).Value2 = 1.048420688269261
$ws.Cells.Item(19, 13).Value2 = 1.057846721700613
$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.036480018410351
$ws.Cells.Item(20, 4).Value2 = 1.043190076482385
$ws.Cells.Item(20, 5).Value2 = 1.044865790938572
$ws.Cells.Item(20, 6).Value2 = 1.054316451929777
$ws.Cells.Item(20, 9).Value2 = 1.036995600007931
$ws.Cells.Item(20, 10).Value2 = 1.042672893042079
$ws.Cells.Item(20, 11).Value2 = 1.046538689028501
$ws.Cells.Item(20, 12).Value2 = 1.048208648022092
$ws.Cells.Item(20, 13).Value2 = 1.057627225836603
$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.035433725455265
$ws.Cells.Item(21, 4).Value2 = 1.042415200477435
$ws.Cells.Item(21, 5).Value2 = 1.043969010472466
$ws.Cells.Item(21, 6).Value2 = 1.053397367928969
$ws.Cells.Item(21, 9).Value2 = 1.036825480093159
$ws.Cells.Item(21, 10).Value2 = 1.042023601040351
$ws.Cells.Item(21, 11).Value2 = 1.045971626247385
$ws.Cells.Item(21, 12).Value2 = 1.047519741394451
$ws.Cells.Item(21, 13).Value2 = 1.056913940838135
$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.034776689082257
$ws.Cells.Item(22, 4).Value2 = 1.041928533477485
$ws.Cells.Item(22, 5).Value2 = 1.043406234796854
$ws.Cells.Item(22, 6).Value2 = 1.052820459154247
$ws.Cells.Item(22, 9).Value2 = 1.036717334381913
$ws.Cells.Item(22, 10).Value2 = 1.041615500887314
$ws.Cells.Item(22, 11).Value2 = 1.045614884474228
$ws.Cells.Item(22, 12).Value2 = 1.047086955511406
$ws.Cells.Item(22, 13).Value2 = 1.056465721045478
$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.035124919314671
$ws.Cells.Item(23, 4).Value2 = 1.04218647421572
$ws.Cells.Item(23, 5).Value2 = 1.04370447123186
$ws.Cells.Item(23, 6).Value2 = 1.053126198206018
$ws.Cells.Item(23, 9).Value2 = 1.036774776545445
$ws.Cells.Item(23, 10).Value2 = 1.041831829245137
$ws.Cells.Item(23, 11).Value2 = 1.045804019227503
$ws.Cells.Item(23, 12).Value2 = 1.047316349199925
$ws.Cells.Item(23, 13).Value2 = 1.056703306573549
$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.036497573275765
$ws.Cells.Item(24, 4).Value2 = 1.043203076211061
$ws.Cells.Item(24, 5).Value2 = 1.044880843522707
$ws.Cells.Item(24, 6).Value2 = 1.054331876564242
$ws.Cells.Item(24, 9).Value2 = 1.036998431964077
$ws.Cells.Item(24, 10).Value2 = 1.042683780714799
$ws.Cells.Item(24, 11).Value2 = 1.046548192310665
$ws.Cells.Item(24, 12).Value2 = 1.048220203586069
$ws.Cells.Item(24, 13).Value2 = 1.05763918829695
$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.038093942032549
$ws.Cells.Item(25, 4).Value2 = 1.044385029244429
$ws.Cells.Item(25, 5).Value2 = 1.04625054119218
$ws.Cells.Item(25, 6).Value2 = 1.055735100230881
$ws.Cells.Item(25, 9).Value2 = 1.037252792765049
$ws.Cells.Item(25, 10).Value2 = 1.043672971631011
$ws.Cells.Item(25, 11).Value2 = 1.047410820412208
$ws.Cells.Item(25, 12).Value2 = 1.049270585566155
$ws.Cells.Item(25, 13).Value2 = 1.058726269477797
